$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Appointments")
$ws.Range("G2").Value = "Lorne NEW TEAM 2"
$ws.Range("H11").Select()
